# Fix IFRS company_list numeric data (BGF.xlsx) — "error solve ifrs list"
# Rewrites the per-row financial metrics in columns D:AJ for rows 2-9,
# and drops the now-unused "U" column cells on rows 7-9 (cell removed,
# not just blanked, to match the source data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 33680
$ws.Range("E2").Value = 1241
$ws.Range("F2").Value = 1241
$ws.Range("G2").Value = 1362
$ws.Range("H2").Value = 1015
$ws.Range("I2").Value = 1016
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = 13385
$ws.Range("L2").Value = 7664
$ws.Range("M2").Value = 5721
$ws.Range("N2").Value = 5590
$ws.Range("O2").Value = 131
$ws.Range("P2").Value = 246
$ws.Range("Q2").Value = 2603
$ws.Range("R2").Value = -2557
$ws.Range("S2").Value = 70
$ws.Range("T2").Value = 678
$ws.Range("U2").Value = 1925
$ws.Range("V2").Value = 1129
$ws.Range("W2").Value = 3.69
$ws.Range("X2").Value = 3.01
$ws.Range("Y2").Value = 27.01
$ws.Range("Z2").Value = 8.119999999999999
$ws.Range("AA2").Value = 133.96
$ws.Range("AB2").Value = 2148.81
$ws.Range("AC2").Value = 2063
$ws.Range("AD2").Value = 19.12
$ws.Range("AE2").Value = 11343
$ws.Range("AF2").Value = 3.48
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 0.76
$ws.Range("AI2").Value = 14.54
$ws.Range("AJ2").Value = 49279859
# Row 3
$ws.Range("D3").Value = 43343
$ws.Range("E3").Value = 1836
$ws.Range("F3").Value = 1836
$ws.Range("G3").Value = 1965
$ws.Range("H3").Value = 1528
$ws.Range("I3").Value = 1517
$ws.Range("J3").Value = 12
$ws.Range("K3").Value = 16216
$ws.Range("L3").Value = 8479
$ws.Range("M3").Value = 7737
$ws.Range("N3").Value = 7732
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 248
$ws.Range("Q3").Value = 3593
$ws.Range("R3").Value = -4446
$ws.Range("S3").Value = 267
$ws.Range("T3").Value = 933
$ws.Range("U3").Value = 2660
$ws.Range("V3").Value = 678
$ws.Range("W3").Value = 4.24
$ws.Range("X3").Value = 3.53
$ws.Range("Y3").Value = 22.77
$ws.Range("Z3").Value = 10.33
$ws.Range("AA3").Value = 109.59
$ws.Range("AB3").Value = 3041.47
$ws.Range("AC3").Value = 3072
$ws.Range("AD3").Value = 28.78
$ws.Range("AE3").Value = 15605
$ws.Range("AF3").Value = 5.66
$ws.Range("AG3").Value = 600
$ws.Range("AH3").Value = 0.68
$ws.Range("AI3").Value = 19.6
$ws.Range("AJ3").Value = 49547625
# Row 4
$ws.Range("D4").Value = 861
$ws.Range("E4").Value = 80
$ws.Range("F4").Value = 2172
$ws.Range("G4").Value = 199
$ws.Range("H4").Value = 1846
$ws.Range("I4").Value = 1835
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 19176
$ws.Range("L4").Value = 9754
$ws.Range("M4").Value = 9422
$ws.Range("N4").Value = 9240
$ws.Range("O4").Value = 182
$ws.Range("P4").Value = 495
$ws.Range("Q4").Value = 3320
$ws.Range("R4").Value = -1580
$ws.Range("S4").Value = -2248
$ws.Range("T4").Value = 1181
$ws.Range("U4").Value = 2139
$ws.Range("V4").Value = 554
$ws.Range("W4").Value = 9.31
$ws.Range("X4").Value = 214.48
$ws.Range("Y4").Value = 21.62
$ws.Range("Z4").Value = 10.43
$ws.Range("AA4").Value = 103.53
$ws.Range("AB4").Value = 1779.47
$ws.Range("AC4").Value = 3704
$ws.Range("AD4").Value = 22.76
$ws.Range("AE4").Value = 18649
$ws.Range("AF4").Value = 4.52
$ws.Range("AG4").Value = 800
$ws.Range("AH4").Value = 0.95
$ws.Range("AI4").Value = 21.6
$ws.Range("AJ4").Value = 49547625
# Row 5
$ws.Range("D5").Value = 1243
$ws.Range("E5").Value = 110
$ws.Range("F5").Value = 110
$ws.Range("G5").Value = -43
$ws.Range("H5").Value = 34979
$ws.Range("I5").Value = 34976
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 8581
$ws.Range("L5").Value = 1329
$ws.Range("M5").Value = 7252
$ws.Range("N5").Value = 7082
$ws.Range("O5").Value = 170
$ws.Range("P5").Value = 323
$ws.Range("Q5").Value = 3841
$ws.Range("R5").Value = -3509
$ws.Range("S5").Value = -1186
$ws.Range("T5").Value = 1746
$ws.Range("U5").Value = 2095
$ws.Range("V5").Value = 243
$ws.Range("W5").Value = 8.81
$ws.Range("X5").Value = 2813.01
$ws.Range("Y5").Value = 428.57
$ws.Range("Z5").Value = 252.03
$ws.Range("AA5").Value = 18.32
$ws.Range("AB5").Value = 13439.63
$ws.Range("AC5").Value = 74962
$ws.Range("AD5").Value = 0.2
$ws.Range("AE5").Value = 21957
$ws.Range("AF5").Value = 0.6899999999999999
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 1.31
$ws.Range("AI5").Value = 0.18
$ws.Range("AJ5").Value = 32263719
# Row 6
$ws.Range("D6").Value = 2206
$ws.Range("E6").Value = 295
$ws.Range("F6").Value = 295
$ws.Range("G6").Value = 523
$ws.Range("H6").Value = 472
$ws.Range("I6").Value = 483
$ws.Range("K6").Value = 16249
$ws.Range("L6").Value = 927
$ws.Range("M6").Value = 15322
$ws.Range("N6").Value = 14991
$ws.Range("P6").Value = 954
$ws.Range("Q6").Value = -292
$ws.Range("R6").Value = 294
$ws.Range("S6").Value = -37
$ws.Range("T6").Value = 40
$ws.Range("U6").Value = -332
$ws.Range("V6").Value = 277
$ws.Range("W6").Value = 13.36
$ws.Range("X6").Value = 21.4
$ws.Range("Y6").Value = 4.37
$ws.Range("Z6").Value = 3.8
$ws.Range("AA6").Value = 6.05
$ws.Range("AB6").Value = 1493.8
$ws.Range("AC6").Value = 575
$ws.Range("AD6").Value = 14.04
$ws.Range("AE6").Value = 15720
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 0.43
$ws.Range("AI6").Value = 6.91
$ws.Range("AJ6").Value = 95369179
# Row 7
$ws.Range("D7").Value = 2322
$ws.Range("E7").Value = 252
$ws.Range("G7").Value = 339
$ws.Range("H7").Value = 288
$ws.Range("I7").Value = 348
$ws.Range("K7").Value = 16750
$ws.Range("L7").Value = 1180
$ws.Range("M7").Value = 15588
$ws.Range("N7").Value = 15308
$ws.Range("P7").Value = 958
$ws.Range("Q7").Value = 2646
$ws.Range("R7").Value = 3607
$ws.Range("S7").Value = -28
$ws.Range("T7").Value = 25
$ws.Range("U7").Value = $null
$ws.Range("W7").Value = 10.83
$ws.Range("X7").Value = 12.41
$ws.Range("Y7").Value = 2.29
$ws.Range("Z7").Value = 1.75
$ws.Range("AA7").Value = 7.57
$ws.Range("AC7").Value = 364
$ws.Range("AD7").Value = 14.11
$ws.Range("AE7").Value = 15994
$ws.Range("AF7").Value = 0.32
$ws.Range("AG7").Value = 35
$ws.Range("AH7").Value = 0.68
$ws.Range("AI7").Value = 9.640000000000001
# Row 8
$ws.Range("D8").Value = 2443
$ws.Range("E8").Value = 360
$ws.Range("G8").Value = 481
$ws.Range("H8").Value = 386
$ws.Range("I8").Value = 445
$ws.Range("K8").Value = 16990
$ws.Range("L8").Value = 1180
$ws.Range("M8").Value = 15937
$ws.Range("N8").Value = 15710
$ws.Range("P8").Value = 958
$ws.Range("Q8").Value = 471
$ws.Range("R8").Value = 3400
$ws.Range("S8").Value = 15
$ws.Range("T8").Value = 20
$ws.Range("U8").Value = $null
$ws.Range("W8").Value = 14.74
$ws.Range("X8").Value = 15.78
$ws.Range("Y8").Value = 2.87
$ws.Range("Z8").Value = 2.28
$ws.Range("AA8").Value = 7.4
$ws.Range("AC8").Value = 465
$ws.Range("AD8").Value = 11.03
$ws.Range("AE8").Value = 16415
$ws.Range("AF8").Value = 0.31
$ws.Range("AG8").Value = 35
$ws.Range("AH8").Value = 0.68
$ws.Range("AI8").Value = 7.53
# Row 9
$ws.Range("D9").Value = 2539
$ws.Range("E9").Value = 398
$ws.Range("G9").Value = 508
$ws.Range("H9").Value = 406
$ws.Range("I9").Value = 471
$ws.Range("K9").Value = 17210
$ws.Range("L9").Value = 1180
$ws.Range("M9").Value = 16304
$ws.Range("N9").Value = 16144
$ws.Range("P9").Value = 958
$ws.Range("Q9").Value = 644
$ws.Range("R9").Value = 3314
$ws.Range("S9").Value = -35
$ws.Range("T9").Value = 20
$ws.Range("U9").Value = $null
$ws.Range("W9").Value = 15.66
$ws.Range("X9").Value = 15.97
$ws.Range("Y9").Value = 2.96
$ws.Range("Z9").Value = 2.37
$ws.Range("AA9").Value = 7.24
$ws.Range("AC9").Value = 492
$ws.Range("AD9").Value = 10.43
$ws.Range("AE9").Value = 16868
$ws.Range("AF9").Value = 0.3
$ws.Range("AG9").Value = 35
$ws.Range("AH9").Value = 0.68
$ws.Range("AI9").Value = 7.11

Write-Output "Updated 241 cells, cleared 3 cells"
